$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh (GitHub Actions scheduled update).
# D-column "Price" cells are plain text (e.g. "25.680.67"), so force
# Text number-format before assigning to stop Excel from re-parsing them
# as numbers (which would also mangle values like "0.000006534").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.682.25'
$ws.Range("E2").Value = '  -3.74%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.745.05'
$ws.Range("E3").Value = '  -5.71%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.69'
$ws.Range("E5").Value = '  -9.95%  '

$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4901'
$ws.Range("E7").Value = '  -8.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.57'
$ws.Range("E8").Value = '  -7.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2486'
$ws.Range("E9").Value = '  -22.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.05959'
$ws.Range("E10").Value = '  -15.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.745.26'
$ws.Range("E11").Value = '  -5.82%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.06782'
$ws.Range("E12").Value = '  -13.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.77'
$ws.Range("E13").Value = '  -22.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.452'
$ws.Range("E14").Value = '  -11.97%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.13'
$ws.Range("E15").Value = '  -13.90%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.5603'
$ws.Range("E16").Value = '  -28.07%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.19%  '

$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.728.59'
$ws.Range("E19").Value = '  -3.62%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.47'
$ws.Range("E20").Value = '  -19.09%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006534'
$ws.Range("E21").Value = '  -18.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.965.14'
$ws.Range("E22").Value = '  -6.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.963'
$ws.Range("E23").Value = '  -14.79%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.002'
$ws.Range("E24").Value = '  -17.38%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.847'
$ws.Range("E25").Value = '  -16.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.26'
$ws.Range("E26").Value = '  -4.72%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.481'
$ws.Range("E27").Value = '  -12.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.809'
$ws.Range("E28").Value = '  -18.55%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '14.61'
$ws.Range("E29").Value = '  -14.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '101.82'
$ws.Range("E30").Value = '  -8.99%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.745'
$ws.Range("E31").Value = '  -12.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08023'
$ws.Range("E32").Value = '  -8.41%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.303'
$ws.Range("E33").Value = '  -19.92%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04392'
$ws.Range("E34").Value = '  -10.09%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9996'
$ws.Range("E35").Value = '  -0.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.582'
$ws.Range("E36").Value = '  -10.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9813'
$ws.Range("E37").Value = '  -14.38%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6024'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.679'
$ws.Range("E39").Value = '  -13.95%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.009'
$ws.Range("E40").Value = '  -14.79%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.35'
$ws.Range("E42").Value = '  -5.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01494'
$ws.Range("E43").Value = '  -14.90%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7571'
$ws.Range("E44").Value = '  -16.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.135'
$ws.Range("E45").Value = '  -13.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.3695'
$ws.Range("E46").Value = '  -23.77%  '

$ws.Range("E47").Value = '  -12.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1066'
$ws.Range("E48").Value = '  -15.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.01'
$ws.Range("E49").Value = '  -14.51%  '

$ws.Range("E50").Value = '  -13.59%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.849'
$ws.Range("E51").Value = '  -24.47%  '
